$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right after the title (Heading1) paragraph.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play Avatar: Gateway Guardians Slot
#    for Free - Review") right before the final paragraph (the one that
#    currently holds the italic "Please create a feature image..."
#    text).
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlSnippet = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Avatar: Gateway Guardians Slot for Free - Review</w:t></w:r></w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertPoint.InsertXML($xmlSnippet) | Out-Null

# InsertXML leaves behind an extra, empty stray paragraph between the
# newly inserted one and the original final paragraph - remove it.
$stray = $d.Paragraphs.Item($n + 1)
$stray.Range.Delete() | Out-Null

# ---------------------------------------------------------------------
# 3. Swap out the old image-prompt text in the final paragraph for the
#    new meta-description text, keeping the existing (italic) run
#    formatting intact.
# ---------------------------------------------------------------------
$oldText = 'Please create a feature image for "Avatar: Gateway Guardians" that fits the following criteria: - It should be in a cartoon style. - The main character in the image should be a happy Maya warrior wearing glasses. The image should feature a round frame, similar to the circular reels in the game. The happy Maya warrior should be standing in the center of the frame wielding a staff adorned with blue and azure hues. The warrior should be wearing a traditional Mayan headdress and gray armor. The background should be misty, with a few floating rocks and a hint of blue and green hues. The overall tone of the image should be vibrant and exciting, inviting players to try their luck in the world of Avatar: Gateway Guardians.'
$newText = 'Explore the innovative features of Avatar: Gateway Guardians slot game and play for free. Read our review for a unique playing experience.'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
